$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newJson = @'
<<ConditionalFormat json={
    "rules": [
        {
            "operator": "between",
            "minValue": 0,
            "maxValue": 1000,
            "background": "#E2FEE2"
        },
        {
            "operator": "between",
            "minValue": 1001,
            "maxValue": 10000,
            "background": "#FFFFD4"
        },
        {
            "operator": "greaterThan",
            "value": 100001,
            "background": "#FB8383"
        }
    ]
}>>
'@

$ws.Range("H4").Value = $newJson
$ws.Range("J4").Value = $newJson

# Setting the (wrapped) cell text above can trigger an row auto-fit in the
# interop layer; restore the original explicit row height so only the text
# itself changes.
$ws.Rows.Item(4).RowHeight = 16.5

$ws.Range("Q11").Select()
